$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.478.54"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.89%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.828.46"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.79%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.13"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.57%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.01%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5083"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -4.41%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3911"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.20%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07640"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.51%  "

$ws.Range("E10").Value = "  +0.94%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.108"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.98%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.11"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.69%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.301"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.02%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.574"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.64%  "

$ws.Range("E15").Value = "  +0.06%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.822.86"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.76%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.30"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +5.57%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001084"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.30%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06672"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.53%  "

$ws.Range("E20").Value = "  +3.08%  "

$ws.Range("E21").Value = "  +0.00%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.166"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.71%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.503.37"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.87%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.15"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.51%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.256"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +7.78%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.76"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.25%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.63"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.50%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.034.68"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.72%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.394"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +4.56%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.08"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.74%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.126"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.19%  "

$ws.Range("E32").Value = "  -0.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.686"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.61%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.660"
$ws.Range("D34").ClearFormats()

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07017"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.12%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2227"
$ws.Range("D36").ClearFormats()

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.994"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +7.09%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02326"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.30%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.139"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.57%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6276"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.62%  "

$ws.Range("E41").Value = "  +0.07%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.184"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.38%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.000"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.03%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.398"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.36%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.48"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.25%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5906"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.51%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.714"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.99%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.51"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.48%  "

$ws.Range("E49").Value = "  +3.66%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.197"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.57%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06922"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.71%  "
